# Rename existing sheet "Feuil1" -> "Details"
$wb = $excel.ActiveWorkbook
$details = $wb.Worksheets.Item(1)
$details.Name = "Details"

# Add a new blank worksheet "Summary" right after "Details"
$summary = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $details)
$summary.Name = "Summary"

# Keep "Details" the active sheet and move its selection from D36 to D41
$details.Activate()
$details.Range("D41").Select() | Out-Null
